$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the lastlogin timestamp for the existing "v" user (row 2)
$ws.Range("E2").Value = "2022-01-10 00:13:10.13S"

# Add a new row (row 4) for the "vikrant" user
$ws.Range("A4").Value = "vikrant"
$ws.Range("B4").Value = "Vikrant"
$ws.Range("C4").Value = "Deshpande"
$ws.Range("D4").Value = "vikrant"
$ws.Range("E4").Value = "2022-01-10 00:12:59.12S"
$ws.Range("F4").Value = 102
